$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for first row (G2)
$wsOverview.Range("G2").Value = "2016-09-03 19:18:36"

# zh-cn sheet: Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2)
$wsZhCn.Range("H2").Value = "2016-09-03 19:18:32"
$wsZhCn.Range("K2").Value = "2016-09-03 19:18:49"

# de-de sheet: Correspond Handoff Datetime (H2) shares the same value as Overview G2,
# and Correspond Handback DateTime (K2)
$wsDeDe.Range("H2").Value = "2016-09-03 19:18:36"
$wsDeDe.Range("K2").Value = "2016-09-03 19:18:57"
